$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing "vada pav" row (row 4), pushing it
# down to row 6, then fill the freed rows 4-5 with the new "login"/"logout"
# entries (siddharth, siddhu).
$ws.Rows("4:5").Insert()

$ws.Range("A4").Value = "siddharth"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1

$ws.Range("A5").Value = "siddhu"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1
